$d = $word.ActiveDocument

$replacements = @(
    @("560×7=3920", "919×8=7352"),
    @("402×9=3618", "506×9=4554"),
    @("458×9=4122", "528×8=4224"),
    @("378×5=1890", "737×9=6633"),
    @("142×5=710",  "254×2=508"),
    @("197×9=1773", "453×8=3624"),
    @("822×8=6576", "423×2=846"),
    @("890×9=8010", "133×3=399"),
    @("984×3=2952", "207×3=621"),
    @("726×5=3630", "123×5=615"),
    @("414×5=2070", "157×5=785"),
    @("701×2=1402", "323×5=1615"),
    @("742×4=2968", "533×3=1599"),
    @("954×5=4770", "826×3=2478"),
    @("143×9=1287", "481×8=3848"),
    @("713×8=5704", "303×4=1212"),
    @("977×8=7816", "325×6=1950"),
    @("767×8=6136", "508×6=3048"),
    @("146×4=584",  "463×6=2778"),
    @("293×9=2637", "165×4=660"),
    @("954×3=2862", "207×4=828"),
    @("533×8=4264", "968×4=3872"),
    @("922×5=4610", "500×7=3500"),
    @("973×3=2919", "357×9=3213"),
    @("555×4=2220", "868×5=4340")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
